$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "66.641.56"
$ws.Range("E2").Value = "  -4.18%  "

$ws.Range("D3").Value = "3.360.31"
$ws.Range("E3").Value = "  -4.86%  "

$ws.Range("E4").Value = "  -0.53%  "

Set-TextValue "D5" "557.06"
$ws.Range("E5").Value = "  -4.01%  "

Set-TextValue "D6" "183.57"
$ws.Range("E6").Value = "  -6.32%  "

Set-TextValue "D7" "0.599"
$ws.Range("E7").Value = "  -2.11%  "

Set-TextValue "D8" "0.999"
$ws.Range("E8").Value = "  -0.18%  "

$ws.Range("D9").Value = "3.353.05"
$ws.Range("E9").Value = "  -4.42%  "

$ws.Range("E10").Value = "  -9.09%  "

$ws.Range("E11").Value = "  -5.20%  "

Set-TextValue "D12" "47.88"
$ws.Range("E12").Value = "  -7.17%  "

$ws.Range("E13").Value = "  -6.08%  "

Set-TextValue "D14" "8.71"
$ws.Range("E14").Value = "  -6.05%  "

$ws.Range("D15").Value = "3.892.36"
$ws.Range("E15").Value = "  -5.88%  "

Set-TextValue "D16" "604.01"
$ws.Range("E16").Value = "  -10.97%  "

$ws.Range("D17").Value = "66.455.49"
$ws.Range("E17").Value = "  -4.80%  "

$ws.Range("D18").Value = "3.351.76"
$ws.Range("E18").Value = "  -5.91%  "

$ws.Range("E19").Value = "  -3.92%  "

$ws.Range("E20").Value = "  -3.46%  "

$ws.Range("E21").Value = "  -5.67%  "

Set-TextValue "D22" "0.911"
$ws.Range("E22").Value = "  -5.63%  "

Set-TextValue "D23" "16.84"
$ws.Range("E23").Value = "  -5.35%  "

Set-TextValue "D24" "5.06"
$ws.Range("E24").Value = "  -1.32%  "

Set-TextValue "D25" "97.94"
$ws.Range("E25").Value = "  -8.74%  "

$ws.Range("E26").Value = "  -7.70%  "

$ws.Range("E27").Value = "  -5.57%  "

$ws.Range("E28").Value = "  -7.32%  "

Set-TextValue "D29" "8.82"
$ws.Range("E29").Value = "  -8.54%  "

Set-TextValue "D30" "30.76"
$ws.Range("E30").Value = "  -7.62%  "

Set-TextValue "D31" "6.35"
$ws.Range("E31").Value = "  -7.26%  "

Set-TextValue "D32" "3.84"
$ws.Range("E32").Value = "  -11.58%  "

Set-TextValue "D33" "11.17"
$ws.Range("E33").Value = "  -5.38%  "

$ws.Range("E34").Value = "  -5.35%  "

$ws.Range("D35").Value = "3.832.77"
$ws.Range("E35").Value = "  +1.74%  "

Set-TextValue "D36" "57.91"
$ws.Range("E36").Value = "  -7.86%  "

Set-TextValue "D37" "526.56"
$ws.Range("E37").Value = "  +5.77%  "

Set-TextValue "D38" "0.998"
$ws.Range("E38").Value = "  +0.04%  "

Set-TextValue "D39" "3.63"
$ws.Range("E39").Value = "  +45.78%  "

$ws.Range("E40").Value = "  -4.99%  "

$ws.Range("D41").Value = "0.0₃0725"
$ws.Range("E41").Value = "  -10.72%  "

$ws.Range("E42").Value = "  -6.54%  "

$ws.Range("E43").Value = "  -5.83%  "

Set-TextValue "D44" "0.350"
$ws.Range("E44").Value = "  -5.86%  "

Set-TextValue "D45" "32.52"
$ws.Range("E45").Value = "  -5.95%  "

Set-TextValue "D46" "0.0418"
$ws.Range("E46").Value = "  -8.50%  "

Set-TextValue "D47" "3.17"
$ws.Range("E47").Value = "  -6.82%  "

$ws.Range("E48").Value = "  -9.19%  "

$ws.Range("E49").Value = "  -4.99%  "

$ws.Range("E50").Value = "  -0.78%  "

Set-TextValue "D51" "7.70"
$ws.Range("E51").Value = "  -7.66%  "
